$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.1614932821685667
$ws.Range("B6").Value = 0.2023608905510596
$ws.Range("B7").Value = 0.2211020455267163
$ws.Range("B8").Value = 0.1766371156411559
$ws.Range("B9").Value = 0.1414627756152198
$ws.Range("B11").Value = 0.226444284693467
$ws.Range("B12").Value = 0.08769308169122851
$ws.Range("B13").Value = 0.1435805489676914
$ws.Range("B14").Value = 0.1624484649590589
$ws.Range("B15").Value = 0.1731355958147772
$ws.Range("B17").Value = 0.1916233172933781
$ws.Range("B19").Value = 0.2287516070550251
$ws.Range("B20").Value = 0.1736546478030555
$ws.Range("B21").Value = 0.1848761074453971
$ws.Range("B22").Value = 0.1713699259368486
$ws.Range("B23").Value = 0.1779696044714125
$ws.Range("B24").Value = 0.1404030286485607
$ws.Range("B25").Value = 0.2022389753616026
$ws.Range("B26").Value = 0.1912384331656097
$ws.Range("B27").Value = 0.2541806191132886
$ws.Range("B29").Value = 0.1019587035736605
$ws.Range("B30").Value = 0.2365928227144809
$ws.Range("B31").Value = 0.2846227786979708
$ws.Range("B32").Value = 0.1766573030530097
$ws.Range("B33").Value = 0.133329554673875
$ws.Range("B34").Value = 0.1900616305511722
$ws.Range("B36").Value = 0.1437153912033442
$ws.Range("B37").Value = 0.2056578078520161
$ws.Range("B38").Value = 0.2345243472556311
$ws.Range("B39").Value = 0.1610943787723742
$ws.Range("B41").Value = 0.1458240755039918
$ws.Range("B42").Value = 0.2059327919499025
$ws.Range("B43").Value = 0.2838018152548218
$ws.Range("B44").Value = 0.1931262228905267
$ws.Range("B46").Value = 0.1199325929446099
$ws.Range("B47").Value = 0.3067553465909482
$ws.Range("B48").Value = 0.1570109909858135
$ws.Range("B49").Value = 0.1468919994367651
$ws.Range("B50").Value = 0.1691606932123232
$ws.Range("B51").Value = 0.1741003377770491
$ws.Range("B52").Value = 0.1707944155372615
$ws.Range("B53").Value = 0.1731564074651663
$ws.Range("B54").Value = 0.2066678131385672
$ws.Range("B56").Value = 0.1307145141641815
$ws.Range("B57").Value = 0.1079885373714472
$ws.Range("B58").Value = 0.19053454548198
$ws.Range("B59").Value = 0.1919971538337066
$ws.Range("B61").Value = 0.2523508369705812
$ws.Range("B62").Value = 0.08738840731918887
$ws.Range("B63").Value = 0.1317842032747991
$ws.Range("B64").Value = 0.2194312291951577
$ws.Range("B65").Value = 0.1060408619463494
$ws.Range("B66").Value = 0.1351086119716099
$ws.Range("B67").Value = 0.1730950473408821
$ws.Range("B68").Value = 0.2893374763815736
$ws.Range("B69").Value = 0.1856236362155468
$ws.Range("B70").Value = 0.2173390525439929
$ws.Range("B72").Value = 0.2793372305797195
$ws.Range("B73").Value = 0.2442622511631537
$ws.Range("B74").Value = 0.1492557098424128
$ws.Range("B75").Value = 0.138045466837423
$ws.Range("B76").Value = 0.1305344264377558
$ws.Range("B77").Value = 0.1447337515325399
$ws.Range("B78").Value = 0.1493412722513493
